$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.262707239097949
$ws.Range("C2").Value = 0.3043599639595413
$ws.Range("E2").Value = 0.6418854245945766
$ws.Range("F2").Value = 2.291294084569103
$ws.Range("G2").Value = 0.5133542502961888
$ws.Range("H2").Value = 0.6295261512474042
$ws.Range("J2").Value = 0.04377196850951037
$ws.Range("B3").Value = 1.120204199138129
$ws.Range("C3").Value = 0.2654254544516164
$ws.Range("E3").Value = 0.6189219973254865
$ws.Range("F3").Value = 2.256770304688501
$ws.Range("G3").Value = 0.51420234392792
$ws.Range("H3").Value = 0.6375345785623807
$ws.Range("J3").Value = 0.04415469535792838
$ws.Range("B4").Value = 1.032743207378019
$ws.Range("C4").Value = 0.2414433423650166
$ws.Range("E4").Value = 0.6050139318468979
$ws.Range("F4").Value = 2.237333643526355
$ws.Range("G4").Value = 0.5155636331247706
$ws.Range("H4").Value = 0.643093673735649
$ws.Range("J4").Value = 0.04445087243192702
$ws.Range("B5").Value = 0.9971116985606727
$ws.Range("C5").Value = 0.2316514638650915
$ws.Range("E5").Value = 0.5993945768389182
$ws.Range("F5").Value = 2.229853652754002
$ws.Range("G5").Value = 0.5163277833372319
$ws.Range("H5").Value = 0.6455196874914435
$ws.Range("J5").Value = 0.04458682901672617
$ws.Range("B6").Value = 0.9911957033572207
$ws.Range("C6").Value = 0.2300243895604126
$ws.Range("E6").Value = 0.5984644078881303
$ws.Range("F6").Value = 2.228638152584836
$ws.Range("G6").Value = 0.5164672608885752
$ws.Range("H6").Value = 0.6459322047868596
$ws.Range("J6").Value = 0.04461032271280274
$ws.Range("B7").Value = 1.032262629146089
$ws.Range("C7").Value = 0.2413113622352228
$ws.Range("E7").Value = 0.6049379515470719
$ws.Range("F7").Value = 2.237230984718266
$ws.Range("G7").Value = 0.5155730934026934
$ws.Range("H7").Value = 0.6431257424209917
$ws.Range("J7").Value = 0.0444526443563511
$ws.Range("B8").Value = 1.213564801243081
$ws.Range("C8").Value = 0.2909512452490617
$ws.Range("E8").Value = 0.6339279370114781
$ws.Range("F8").Value = 2.279023409295348
$ws.Range("G8").Value = 0.5134710397386044
$ws.Range("H8").Value = 0.632153768914435
$ws.Range("J8").Value = 0.04389116071777366
$ws.Range("B9").Value = 1.569390280894027
$ws.Range("C9").Value = 0.3876880548206145
$ws.Range("E9").Value = 0.6922956342516784
$ws.Range("F9").Value = 2.375060163019526
$ws.Range("G9").Value = 0.5161069903113003
$ws.Range("H9").Value = 0.61576573922477
$ws.Range("J9").Value = 0.04328109603355301
$ws.Range("B10").Value = 1.831028095596139
$ws.Range("C10").Value = 0.4583952198727275
$ws.Range("E10").Value = 0.7361074719727014
$ws.Range("F10").Value = 2.454364863858245
$ws.Range("G10").Value = 0.5222871558708135
$ws.Range("H10").Value = 0.6069000738564796
$ws.Range("J10").Value = 0.04313995097061962
$ws.Range("B11").Value = 1.950110164522357
$ws.Range("C11").Value = 0.4904841077038213
$ws.Range("E11").Value = 0.7562415852600566
$ws.Range("F11").Value = 2.492377713296236
$ws.Range("G11").Value = 0.526046851192703
$ws.Range("H11").Value = 0.6035666958547807
$ws.Range("J11").Value = 0.04314409573635203
$ws.Range("B12").Value = 1.99521260407613
$ws.Range("C12").Value = 0.5026243824529502
$ws.Range("E12").Value = 0.7638951721974223
$ws.Range("F12").Value = 2.507053417158573
$ws.Range("G12").Value = 0.5276090767728761
$ws.Range("H12").Value = 0.6024059053908672
$ws.Range("J12").Value = 0.04315563113028986
$ws.Range("B13").Value = 1.985498605381451
$ws.Range("C13").Value = 0.5000102523916894
$ws.Range("E13").Value = 0.7622455362973
$ws.Range("F13").Value = 2.503880200330514
$ws.Range("G13").Value = 0.5272664290009601
$ws.Range("H13").Value = 0.6026513746796098
$ws.Range("J13").Value = 0.04315270135444038
$ws.Range("B14").Value = 1.953820597255117
$ws.Range("C14").Value = 0.4914831184315176
$ws.Range("E14").Value = 0.7568706649914816
$ws.Range("F14").Value = 2.493579444602773
$ws.Range("G14").Value = 0.5261725881574932
$ws.Range("H14").Value = 0.6034691579226745
$ws.Range("J14").Value = 0.04314484424232035
$ws.Range("B15").Value = 1.934418022410568
$ws.Range("C15").Value = 0.4862585544527747
$ws.Range("E15").Value = 0.7535822060880832
$ws.Range("F15").Value = 2.487306613433987
$ws.Range("G15").Value = 0.5255206810236928
$ws.Range("H15").Value = 0.6039833165542206
$ws.Range("J15").Value = 0.04314133345982896
$ws.Range("B16").Value = 1.8232470233595
$ws.Range("C16").Value = 0.4562965908646675
$ws.Range("E16").Value = 0.7347957557001337
$ws.Range("F16").Value = 2.451919825331231
$ws.Range("G16").Value = 0.5220607164939537
$ws.Range("H16").Value = 0.6071320717572632
$ws.Range("J16").Value = 0.0431410672073298
$ws.Range("B17").Value = 1.755062842977509
$ws.Range("C17").Value = 0.4378962872710304
$ws.Range("E17").Value = 0.7233230350062314
$ws.Range("F17").Value = 2.430709037585387
$ws.Range("G17").Value = 0.5201825483140965
$ws.Range("H17").Value = 0.6092435531190574
$ws.Range("J17").Value = 0.04315850769099683
$ws.Range("B18").Value = 1.715850925208372
$ws.Range("C18").Value = 0.427305758519708
$ws.Range("E18").Value = 0.7167434248640632
$ws.Range("F18").Value = 2.41869129423182
$ws.Range("G18").Value = 0.5191914061112755
$ws.Range("H18").Value = 0.610523823729082
$ws.Range("J18").Value = 0.04317496441569091
$ws.Range("B19").Value = 1.702575450945687
$ws.Range("C19").Value = 0.4237187620797158
$ws.Range("E19").Value = 0.7145189814543613
$ws.Range("F19").Value = 2.414653498039286
$ws.Range("G19").Value = 0.5188710653902433
$ws.Range("H19").Value = 0.6109685763073145
$ws.Range("J19").Value = 0.04318163579388568
$ws.Range("B20").Value = 1.76232056490727
$ws.Range("C20").Value = 0.4398557731594224
$ws.Range("E20").Value = 0.724542339695617
$ws.Range("F20").Value = 2.43294809284825
$ws.Range("G20").Value = 0.5203732423661762
$ws.Range("H20").Value = 0.6090119657098683
$ws.Range("J20").Value = 0.0431559851263934
$ws.Range("B21").Value = 1.963124963097926
$ws.Range("C21").Value = 0.4939880464108342
$ws.Range("E21").Value = 0.7584486021672205
$ws.Range("F21").Value = 2.496597377744251
$ws.Range("G21").Value = 0.5264900994547901
$ws.Range("H21").Value = 0.6032261938040477
$ws.Range("J21").Value = 0.04314688050128979
$ws.Range("B22").Value = 2.094412561777688
$ws.Range("C22").Value = 0.5293019173515177
$ws.Range("E22").Value = 0.780778677945392
$ws.Range("F22").Value = 2.539835321723388
$ws.Range("G22").Value = 0.5312961661383326
$ws.Range("H22").Value = 0.6000368144758994
$ws.Range("J22").Value = 0.04319908286653984
$ws.Range("B23").Value = 2.02433735825565
$ws.Range("C23").Value = 0.5104602019881099
$ws.Range("E23").Value = 0.7688451322318741
$ws.Range("F23").Value = 2.51660755351989
$ws.Range("G23").Value = 0.5286564031544145
$ws.Range("H23").Value = 0.6016845909596782
$ws.Range("J23").Value = 0.04316585506472492
$ws.Range("B24").Value = 1.75903938595286
$ws.Range("C24").Value = 0.4389699268036793
$ws.Range("E24").Value = 0.7239910415514146
$ws.Range("F24").Value = 2.431935266000465
$ws.Range("G24").Value = 0.5202867536387004
$ws.Range("H24").Value = 0.6091164597605854
$ws.Range("J24").Value = 0.04315710556034347
$ws.Range("B25").Value = 1.473095370402177
$ws.Range("C25").Value = 0.3615827341006366
$ws.Range("E25").Value = 0.6763427491975307
$ws.Range("F25").Value = 2.347555472990592
$ws.Range("G25").Value = 0.5146576417684798
$ws.Range("H25").Value = 0.6196451543245018
$ws.Range("J25").Value = 0.04339277174229039
